$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The DynamiteNavigation column value for the "News 1" news item is changed
# from "News" to "HR" (a news item cannot be tagged with "news" itself, it
# must be tagged with an actual sub-category, here "HR").
$ws.Range("E2").Value = "HR"

# Update the active selection to reflect the edited cell.
$ws.Range("E2").Select()
